# Applies the "row 27 data entry + new conditional formatting + selection"
# edit described by the commit "AutoCommit_14 июня 2024 г. 13:56:12_SibNout2023".
#
# Summary of the change:
#   - Row 27 (student #24): homework columns C/D get regraded to 5 (green,
#     thick-left/right-border style), E/F are no longer used and are cleared,
#     H gets a grade of 5, a new lab-column entry I27=5 is recorded, and a new
#     column M27=3 is added (extends the used range out to column M).
#   - A new color-scale conditional format is added over E27 (mirrors the
#     existing per-cell color scales on D4/D11/D14).
#   - The active selection is left on D27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 27 cell edits
# ---------------------------------------------------------------------

# E27 and F27 are no longer graded for this student -> remove them entirely.
$ws.Range("E27:F27").Clear()

# C27/D27 move from the "2" grading style to the "5" (fully graded) style.
# Copy formats from a cell that already uses that exact style so no new
# cell style gets created, then overwrite the values.
$ws.Range("I9").Copy()
$ws.Range("C27:D27").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 5

# I27 is a new entry (lab grade) using the same style already used by I16.
$ws.Range("I16").Copy()
$ws.Range("I27").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("I27").Value = 5

# H27 gets graded.
$ws.Range("H27").Value = 5

# M27 is a brand new column entry for this row.
$ws.Range("M27").Value = 3

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# New conditional formatting rule on E27 (same 3-color scale used
# elsewhere on the sheet for D4 / D11 / D14).
# ---------------------------------------------------------------------
$tempRule = $ws.Range("E27").FormatConditions.AddColorScale(3)
$tempRule.SetFirstPriority()
$newRule = $ws.Range("E27").FormatConditions.AddColorScale(3)
$newRule.SetFirstPriority()
$tempRule.Delete()

# ---------------------------------------------------------------------
# Leave the active cell / selection on D27.
# ---------------------------------------------------------------------
$null = $ws.Range("D27").Select()
